# "Regenerate merged AHB files"
#
# The header row (row 1, columns A:U) carries the column captions for a
# side-by-side "AHB diff" table: the first 10 columns describe the "old"
# formulars version, column 11 is the literal "diff" marker, and the last
# 10 columns describe the "new" formulars version. This revision renames
# those generic _old/_new suffixes to the concrete formular versions being
# compared (FV2304 / FV2310), wraps the whole A1:U57 range in a real Excel
# Table (so the generated workbook is filterable/sortable), and freezes
# the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header captions -------------------------------------------------
# Columns A:J ( 1..10) -> "<Label>_old" => "<Label>_FV2304"
# Column  K   (   11 ) -> "diff"        (unchanged)
# Columns L:U (12..21) -> "<Label>_new" => "<Label>_FV2310"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value().Replace("_old", "_FV2304")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value().Replace("_new", "_FV2310")
}

# --- 2. Turn the data range into an Excel Table ---------------------------------
$dataRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row so it stays put while scrolling -------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $null
